$wb = $excel.ActiveWorkbook

# --- Update timestamps on the "data" sheet (column F) ---
$ws = $wb.Worksheets.Item("data")
$ws.Range("F2").Value = "2021-10-05 14:35:12.262562"
$ws.Range("F3").Value = "2021-10-05 14:35:12.262569"
$ws.Range("F4").Value = "2021-10-05 14:35:12.262573"
$ws.Range("F5").Value = "2021-10-05 14:35:12.262575"
$ws.Range("F6").Value = "2021-10-05 14:35:12.262578"
$ws.Range("F7").Value = "2021-10-05 14:35:12.262581"
$ws.Range("F8").Value = "2021-10-05 14:35:12.262584"
$ws.Range("F9").Value = "2021-10-05 14:35:12.262586"
$ws.Range("F10").Value = "2021-10-05 14:35:12.262589"
$ws.Range("F11").Value = "2021-10-05 14:35:12.262592"
$ws.Range("F12").Value = "2021-10-05 14:35:12.262594"

# --- Add the "metadata" sheet (placed after "data") ---
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"
$ws.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$meta.Range("G1").PasteSpecial(-4122)

$meta.Range("A2").Value = 0
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Range("B2").Value = "Periventricular Grey Matter Heterotopia"
$meta.Range("C2").Value = 19
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.0"
$meta.Range("D2").ClearFormats()
$meta.Range("E2").Value = "2021-09-05T03:19:07.209985Z"
$meta.Range("F2").Value = "2021-10-05 14:35:12.258748"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/19/?format=json"

# Keep "data" as the active sheet (matches original workbook view)
$ws.Activate()
